$wb = $excel.ActiveWorkbook

# Sheet "Overview": update the status for the ffa9b05e-... file (row 3)
# from "Ready for handoff" to "Handed back: in sync with en-US" for both
# the zh-cn and de-de columns.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C3").Value = "Handed back: in sync with en-US"

# Sheet "zh-cn": the ffa9b05e-... row (row 3) has been handed back.
# Update its Status and stamp the Latest Handback DateTime.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("H3").Value = "2016-03-11 14:33:21"

# Sheet "de-de": same update for the ffa9b05e-... row (row 3).
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("H3").Value = "2016-03-11 14:33:26"
